$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '64.301.38'
Set-TextValue "E2" '  -6.09%  '
Set-TextValue "D3" '3.357.11'
Set-TextValue "E3" '  -7.42%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.26%  '
Set-TextValue "D5" '182.88'
Set-TextValue "E5" '  -9.36%  '
Set-TextValue "D6" '522.66'
Set-TextValue "E6" '  -10.21%  '
Set-TextValue "D7" '0.595'
Set-TextValue "E7" '  -4.44%  '
Set-TextValue "D8" '3.351.82'
Set-TextValue "E8" '  -7.34%  '
Set-TextValue "E9" '  +0.07%  '
Set-TextValue "D10" '0.611'
Set-TextValue "E10" '  -11.00%  '
Set-TextValue "D11" '56.52'
Set-TextValue "E11" '  -7.30%  '
Set-TextValue "D12" '0.130'
Set-TextValue "E12" '  -13.61%  '
Set-TextValue "D13" '0.0000248'
Set-TextValue "E13" '  -13.21%  '
Set-TextValue "D14" '9.06'
Set-TextValue "E14" '  -10.86%  '
Set-TextValue "D15" '3.909.49'
Set-TextValue "E15" '  -6.97%  '
Set-TextValue "E16" '  -4.60%  '
Set-TextValue "D17" '3.371.90'
Set-TextValue "E17" '  -6.98%  '
Set-TextValue "D18" '64.264.99'
Set-TextValue "E18" '  -5.85%  '
Set-TextValue "D19" '17.11'
Set-TextValue "E19" '  -11.59%  '
Set-TextValue "D20" '10.86'
Set-TextValue "E20" '  -13.22%  '
Set-TextValue "D21" '0.950'
Set-TextValue "E21" '  -11.93%  '
Set-TextValue "D22" '367.44'
Set-TextValue "E22" '  -9.71%  '
Set-TextValue "D23" '79.97'
Set-TextValue "E23" '  -6.84%  '
Set-TextValue "D24" '3.65'
Set-TextValue "E24" '  -14.64%  '
Set-TextValue "D25" '10.55'
Set-TextValue "E25" '  -19.48%  '
Set-TextValue "B26" 'LEO'
Set-TextValue "C26" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D26" '5.88'
Set-TextValue "E26" '  -4.26%  '
Set-TextValue "B27" 'Toncoin'
Set-TextValue "C27" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D27" '3.68'
Set-TextValue "E27" '  -7.97%  '
Set-TextValue "D28" '2.60'
Set-TextValue "E28" '  -11.61%  '
Set-TextValue "D29" '11.11'
Set-TextValue "E29" '  -12.39%  '
Set-TextValue "D30" '8.22'
Set-TextValue "E30" '  -12.98%  '
Set-TextValue "D31" '668.38'
Set-TextValue "E31" '  -1.99%  '
Set-TextValue "D32" '28.59'
Set-TextValue "E32" '  -10.30%  '
Set-TextValue "D33" '6.61'
Set-TextValue "E33" '  -16.02%  '
Set-TextValue "D34" '10.93'
Set-TextValue "E34" '  -11.33%  '
Set-TextValue "D35" '59.37'
Set-TextValue "E35" '  -7.51%  '
Set-TextValue "D36" '0.102'
Set-TextValue "E36" '  -11.17%  '
Set-TextValue "D37" '0.999'
Set-TextValue "E37" '  -0.10%  '
Set-TextValue "D38" '35.71'
Set-TextValue "E38" '  -15.22%  '
Set-TextValue "D39" '0.369'
Set-TextValue "E39" '  -11.76%  '
Set-TextValue "D40" '1.00'
Set-TextValue "E40" '  +0.30%  '
Set-TextValue "D41" '0.125'
Set-TextValue "E41" '  -7.88%  '
Set-TextValue "B42" 'ThetaToken'
Set-TextValue "C42" 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue "D42" '2.71'
Set-TextValue "E42" '  -15.74%  '
Set-TextValue "B43" 'Maker'
Set-TextValue "C43" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D43" '2.762.77'
Set-TextValue "E43" '  -13.69%  '
Set-TextValue "D44" '2.59'
Set-TextValue "E44" '  -9.49%  '
Set-TextValue "B45" 'VeChain'
Set-TextValue "C45" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D45" '0.0380'
Set-TextValue "E45" '  -9.28%  '
Set-TextValue "B46" 'PEPE'
Set-TextValue "C46" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D46" '0.0₃0603'
Set-TextValue "E46" '  -21.99%  '
Set-TextValue "D47" '2.26'
Set-TextValue "E47" '  -16.96%  '
Set-TextValue "D48" '0.123'
Set-TextValue "E48" '  -6.93%  '
Set-TextValue "D49" '134.49'
Set-TextValue "E49" '  -2.90%  '
Set-TextValue "D50" '2.83'
Set-TextValue "E50" '  -7.51%  '
Set-TextValue "D51" '2.54'
Set-TextValue "E51" '  -7.26%  '

Write-Host "Applied cryptos list update"
